$d = $word.ActiveDocument

# --- Step 1: Remove the old "Created fraud detection systems for campaign
#     finance data analysis across multi-terabyte datasets" bullet paragraph.
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*Created fraud detection systems for campaign finance data analysis across multi-terabyte datasets*") {
        $d.Paragraphs($i).Range.Delete()
    }
}

# --- Step 2: Insert three new bullet paragraphs immediately before the
#     "Developed and deployed custom analytical tools..." bullet paragraph.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*Developed and deployed custom analytical tools and algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering*") {
        $anchor = $d.Paragraphs($i).Range
        $anchor.InsertParagraphBefore()
        $anchor.InsertParagraphBefore()
        $anchor.InsertParagraphBefore()

        $p1 = $d.Paragraphs($i).Range
        $p1.Text = "• Developed meta-analytical techniques that identified systematic data quality issues across 20+ years of voter registration data"

        $p2 = $d.Paragraphs($i + 1).Range
        $p2.Text = "• Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters"

        $p3 = $d.Paragraphs($i + 2).Range
        $p3.Text = "• Created fraud detection systems analyzing 5+ terabyte datasets, uncovering demographic miscoding patterns across 2,000+ precincts"

        break
    }
}
